$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 5: the table's style changes from the default "Table_0" style
#    ({1CB131BB-C6DA-449F-B97D-CBA54E4A9D94}) to the built-in table style
#    {B6B57CD2-6733-481B-8DA8-11F7D58745B8}.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shape = $tableSlide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{B6B57CD2-6733-481B-8DA8-11F7D58745B8}", $true)
    }
}

# ---------------------------------------------------------------------------
# 2) The presentation's main theme (used by the slide master / all slides)
#    switches from the "Integral" / "Red Violet" palette to the standard
#    "Office Theme" / "Office" palette.
#    Helper turns an R,G,B triple into the BGR-packed integer PowerPoint's
#    RGB color properties expect.
# ---------------------------------------------------------------------------
function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1
$tcs.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1
$tcs.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$tcs.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$tcs.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$tcs.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$tcs.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$tcs.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$tcs.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$tcs.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$tcs.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$tcs.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
